# 4.0.3 model and data
#
# The "Boolean" sheet lists InputData pathnames of CSV files whose values are
# constrained to boolean data types. Two of the rows - the aggregate
# "trans/BVTQaZ/BVTQaZ.csv" and "trans/VTQaZ/VTQaZ.csv" pathnames - are split
# out into six per-mode files apiece (LDVs, HDVs, aircraft, rail, ships,
# motorbikes).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Boolean")

# --- Expand "trans/BVTQaZ/BVTQaZ.csv" (currently row 17) into six rows ---
# Insert 5 additional blank rows above row 17 so the original row plus the
# five new ones can hold all six new pathnames.
for ($i = 0; $i -lt 5; $i++) {
    $ws.Rows.Item(17).Insert()
}

$ws.Cells.Item(17, 1).Value = "trans/BVTQaZ/BVTQaZ-LDVs.csv"
$ws.Cells.Item(18, 1).Value = "trans/BVTQaZ/BVTQaZ-HDVs.csv"
$ws.Cells.Item(19, 1).Value = "trans/BVTQaZ/BVTQaZ-aircraft.csv"
$ws.Cells.Item(20, 1).Value = "trans/BVTQaZ/BVTQaZ-rail.csv"
$ws.Cells.Item(21, 1).Value = "trans/BVTQaZ/BVTQaZ-ships.csv"
$ws.Cells.Item(22, 1).Value = "trans/BVTQaZ/BVTQaZ-motorbikes.csv"

# --- Expand "trans/VTQaZ/VTQaZ.csv" (now shifted down to row 26) into six rows ---
for ($i = 0; $i -lt 5; $i++) {
    $ws.Rows.Item(26).Insert()
}

$ws.Cells.Item(26, 1).Value = "trans/VTQaZ/VTQaZ-LDVs.csv"
$ws.Cells.Item(27, 1).Value = "trans/VTQaZ/VTQaZ-HDVs.csv"
$ws.Cells.Item(28, 1).Value = "trans/VTQaZ/VTQaZ-aircraft.csv"
$ws.Cells.Item(29, 1).Value = "trans/VTQaZ/VTQaZ-rail.csv"
$ws.Cells.Item(30, 1).Value = "trans/VTQaZ/VTQaZ-ships.csv"
$ws.Cells.Item(31, 1).Value = "trans/VTQaZ/VTQaZ-motorbikes.csv"

# Leave the cursor resting on the last new row, matching where the editor's
# selection ended up on this sheet.
$ws.Activate()
[void]$ws.Range("A32").Select()

# The "Integer" sheet's selection was left on A13 (below its data).
$integer = $wb.Worksheets.Item("Integer")
$integer.Activate()
[void]$integer.Range("A13").Select()

# The "About" sheet was the active tab when the workbook was last saved.
$about = $wb.Worksheets.Item("About")
$about.Activate()
